$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 3670.81
$ws.Range("B32").Value = 12199.68
$ws.Range("F61").Value = 47
$ws.Range("G61").Value = 3301.28
$ws.Range("F70").Value = 1
$ws.Range("G70").Value = 134.95
$ws.Range("F75").Value = 2
$ws.Range("G75").Value = 740.36
$ws.Range("F77").Value = 232
$ws.Range("G77").Value = 10843.68
$ws.Range("F79").Value = 70
$ws.Range("G79").Value = 4351.9
$ws.Range("F80").Value = 6
$ws.Range("G80").Value = 1476.42
$ws.Range("F83").Value = 99
$ws.Range("G83").Value = 14916.33
$ws.Range("B90").Value = 163258.18
$ws.Range("F102").Value = 2
$ws.Range("G102").Value = 98.95999999999999
$ws.Range("B104").Value = 66.94
$ws.Range("F115").Value = 175
$ws.Range("G115").Value = 16941.75
$ws.Range("B117").Value = 10782.81
$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 120.69
$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45
$ws.Range("F151").Value = 84
$ws.Range("G151").Value = 7297.92
$ws.Range("B156").Value = 28147.59
$ws.Range("F160").Value = 11
$ws.Range("G160").Value = 1055.01
$ws.Range("B161").Value = 1094.76
$ws.Range("F183").Value = 1
$ws.Range("G183").Value = 139.35
$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2
$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("B216").Value = 30505.18
$ws.Range("B229").Value = 57802
$ws.Range("E229").Value = 162.71
$ws.Range("F229").Value = -79
$ws.Range("G229").Value = -11334.92
$ws.Range("B230").Value = 63531
$ws.Range("E230").Value = 152.53
$ws.Range("F230").Value = 50
$ws.Range("G230").Value = 7174
$ws.Range("F234").Value = 35
$ws.Range("G234").Value = 1796.2
$ws.Range("F237").Value = 3
$ws.Range("G237").Value = 908.37
$ws.Range("F255").Value = 506
$ws.Range("G255").Value = 86692.98
$ws.Range("B260").Value = 163878.37
$ws.Range("F290").Value = 0
$ws.Range("G290").Value = 0
$ws.Range("F302").Value = 26
$ws.Range("G302").Value = 5483.14
$ws.Range("F303").Value = 17
$ws.Range("G303").Value = 3585.13
$ws.Range("B304").Value = 159749.43
$ws.Range("F320").Value = 34
$ws.Range("G320").Value = 2334.1
$ws.Range("F321").Value = 43
$ws.Range("G321").Value = 2361.56
$ws.Range("F328").Value = 31
$ws.Range("G328").Value = 1153.51
$ws.Range("B330").Value = 24899.12
$ws.Range("F338").Value = 69
$ws.Range("G338").Value = 1635.3
$ws.Range("F343").Value = 28
$ws.Range("G343").Value = 2015.16
$ws.Range("F345").Value = 29
$ws.Range("G345").Value = 1780.89
$ws.Range("B346").Value = 23082
$ws.Range("F350").Value = 61
$ws.Range("G350").Value = 4680.53
$ws.Range("F355").Value = 12
$ws.Range("G355").Value = 1961.4
$ws.Range("F357").Value = 4
$ws.Range("G357").Value = 1045.2
$ws.Range("B358").Value = 33877.51
$ws.Range("B364").Value = 65068
$ws.Range("E364").Value = 13.97
$ws.Range("F364").Value = 63
$ws.Range("G364").Value = 828.45
$ws.Range("B365").Value = 53602
$ws.Range("E365").Value = 15.69
$ws.Range("F365").Value = -231
$ws.Range("G365").Value = -3037.65
$ws.Range("F408").Value = 0
$ws.Range("G408").Value = 0
$ws.Range("F409").Value = 8
$ws.Range("G409").Value = 4671.6
$ws.Range("B411").Value = 4671.6
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("F444").Value = 39
$ws.Range("G444").Value = 2124.72
$ws.Range("B445").Value = 6998.79
$ws.Range("B463").Value = 60025
$ws.Range("E463").Value = 37.22
$ws.Range("F463").Value = -98
$ws.Range("G463").Value = -3217.34
$ws.Range("B464").Value = 64833
$ws.Range("E464").Value = 34.9
$ws.Range("F464").Value = 95
$ws.Range("G464").Value = 3118.85
$ws.Range("F508").Value = 51
$ws.Range("G508").Value = 5300.94
$ws.Range("F509").Value = 188
$ws.Range("G509").Value = 15111.44
$ws.Range("B510").Value = 20412.38
$ws.Range("F534").Value = 1
$ws.Range("G534").Value = 319.52
$ws.Range("F539").Value = 21
$ws.Range("G539").Value = 5440.47
$ws.Range("F540").Value = 0
$ws.Range("G540").Value = 0
$ws.Range("F541").Value = 0
$ws.Range("G541").Value = 0
$ws.Range("F542").Value = 44
$ws.Range("G542").Value = 5699.32
$ws.Range("F544").Value = 0
$ws.Range("G544").Value = 0
$ws.Range("F545").Value = 0
$ws.Range("G545").Value = 0
$ws.Range("F546").Value = 0
$ws.Range("G546").Value = 0
$ws.Range("B547").Value = 13171.39
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 18
$ws.Range("G572").Value = 735.66
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F577").Value = 36
$ws.Range("G577").Value = 1547.64
$ws.Range("B583").Value = 12154.79
$ws.Range("F589").Value = 0
$ws.Range("G589").Value = 0
$ws.Range("B593").Value = 2013.85
$ws.Range("F599").Value = 1248
$ws.Range("G599").Value = 203561.28
$ws.Range("F601").Value = 352
$ws.Range("G601").Value = 99570.24000000001
$ws.Range("F602").Value = 304
$ws.Range("G602").Value = 43973.6
$ws.Range("B606").Value = 347953.17
$ws.Range("B619").Value = 1524995.34
$ws.Range("B620").Value = 1524995.34
